$d = $word.ActiveDocument

function ReplaceText([string]$old, [string]$new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $old"
    }
    return $ok
}

# 1. "560" followed by a _GoBack bookmark that should be dropped; merge the
#    two adjoining runs back together via Find/Replace (Word naturally
#    re-writes the matched range as a single run, and the bookmark pair
#    that sat between the old runs is not part of the Range so it is
#    dropped along with the old run split).
ReplaceText "560W or above" "560W or above"

# 2. "[iw] kW Hybrid Inverter" -> "[iw] Hybrid Inverter"
ReplaceText " kW Hybrid Inverter" " Hybrid Inverter"

# 3. "Daewoo Deep cycle" -> "[bn]"
ReplaceText "Daewoo Deep cycle" "[bn]"

# 4. "180 AH& 12 Vdc" -> "[bs]"
ReplaceText "180 AH& 12 Vdc" "[bs]"

# 5. Greevo (Pvt.) Limited paragraph - two separate runs (bold name + body)
ReplaceText "Greevo (Pvt.) Limited. " "Greevo (Pvt.) Limited. "
ReplaceText "Order by Customer will be deemed confirmed/accepted upon receipt of advance payment as per Payment Option selected by the Customer. Final Price of System will be calculated on USD-PKR exchange rate as at the date of installation or date of final payment by the customer. Balance, if any, due on account of exchange rate variation shall be paid by Customer upon invoice raised by Greevo (Pvt." "Order by Customer will be deemed confirmed/accepted upon receipt of advance payment as per Payment Option selected by the Customer. Final Price of System will be calculated on USD-PKR exchange rate as at the date of installation or date of final payment by the customer. Balance, if any, due on account of exchange rate variation shall be paid by Customer upon invoice raised by Greevo (Pvt."

# 6. "In case of delayed or unpaid payments..."
ReplaceText "In case of delayed or unpaid payments, after-sales services will not be provided, and M/S Greevo Pvt Ltd reserves the right to dismantle and seize system components equal to the outstanding amount." "In case of delayed or unpaid payments, after-sales services will not be provided, and M/S Greevo Pvt Ltd reserves the right to dismantle and seize system components equal to the outstanding amount."

# 7. "Any modifications or alterations..."
ReplaceText "Any modifications or alterations to the system without prior consent from Greevo will nullify after-sales support." "Any modifications or alterations to the system without prior consent from Greevo will nullify after-sales support."

# 8. "Greevo will assist Customer in procuring Net Metering...."
ReplaceText "Greevo will assist Customer in procuring Net Metering. The obligation to procure required approvals rests with the Customer. In no circumstances shall Greevo be held liable on any account" "Greevo will assist Customer in procuring Net Metering. The obligation to procure required approvals rests with the Customer. In no circumstances shall Greevo be held liable on any account"

# 9. "If the said approval is not granted by NEPRA. ..."
ReplaceText "If the said approval is not granted by NEPRA. Greevo will assist Customer to prepare the required documentation for Net-Metering." "If the said approval is not granted by NEPRA. Greevo will assist Customer to prepare the required documentation for Net-Metering."

# 10. Signature block "Greevo PVT Ltd"
ReplaceText "Greevo PVT Ltd" "Greevo PVT Ltd"

Write-Host "Done"
